$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "<one>"
$ws.Range("B2").Value = "<one>"
$ws.Range("C2").Value = 14

# Row 3
$ws.Range("A3").Value = "<oil>"
$ws.Range("B3").Value = "<oil>"

# Row 4
$ws.Range("A4").Value = "<up>"
$ws.Range("B4").Value = "<up>"
$ws.Range("C4").Value = 12

# Row 5
$ws.Range("A5").Value = "<up>"
$ws.Range("B5").Value = "<up>"
$ws.Range("C5").Value = 16

# Row 6
$ws.Range("A6").Value = "<like>"
$ws.Range("B6").Value = "<like>"
$ws.Range("C6").Value = 12

# Row 7
$ws.Range("A7").Value = "<and>"
$ws.Range("B7").Value = "<an>"
$ws.Range("C7").Value = 20

# Row 8
$ws.Range("A8").Value = "<shift>"
$ws.Range("B8").Value = "<shift>"
$ws.Range("C8").Value = 15

# Row 9
$ws.Range("A9").Value = "<said>"
$ws.Range("B9").Value = "<said>"
$ws.Range("C9").Value = 11

# Row 10
$ws.Range("A10").Value = "<tango>"
$ws.Range("B10").Value = "<hine>"
$ws.Range("C10").Value = 14

# Row 11
$ws.Range("A11").Value = "<mike>"
$ws.Range("B11").Value = "<mike>"
$ws.Range("C11").Value = 15

# Row 12
$ws.Range("A12").Value = "<kilo>"
$ws.Range("B12").Value = "<kilo>"
$ws.Range("C12").Value = 10

# Row 13
$ws.Range("A13").Value = "<a>"
$ws.Range("B13").Value = "<a>"
$ws.Range("C13").Value = 17

# Row 14
$ws.Range("A14").Value = "<may>"
$ws.Range("B14").Value = "<may>"
$ws.Range("C14").Value = 12

# Row 15
$ws.Range("A15").Value = "<water>"
$ws.Range("B15").Value = "<water>"
$ws.Range("C15").Value = 15

# Row 16
$ws.Range("A16").Value = "<zulu>"
$ws.Range("B16").Value = "<zulu>"
$ws.Range("C16").Value = 10

# Row 17
$ws.Range("A17").Value = "<water>"
$ws.Range("B17").Value = "<water>"
$ws.Range("C17").Value = 21

# Row 18
$ws.Range("A18").Value = "<uniform>"
$ws.Range("B18").Value = "<in>"
$ws.Range("C18").Value = 16
